# Apply edits described by the diff:
# - Add AVERAGE formulas in column I (rows 4-11)
# - Add IF(...,"GEÇTİ","KALDI") formulas in column J (rows 4-11)
# - Add SUM totals in row 12 (D:G, I)
# - Add AVERAGE totals in row 13 (D:G) and SUM in I13
# - Fill in student number (F17), name (F18) and department (F19)
# - Update the active selection to F19:H19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I: average of the four grade columns (D:G) for each student row
$ws.Range("I4").Formula = "=AVERAGE(D4:G4)"
$ws.Range("I5:I11").Formula = "=AVERAGE(D5:G5)"

# Column J: pass/fail result based on the average in column I
$ws.Range("J4").Formula = '=IF(I4>49, "GEÇTİ", "KALDI")'
$ws.Range("J5:J11").Formula = '=IF(I5>49, "GEÇTİ", "KALDI")'

# Row 12 ("Toplam"): column totals
$ws.Range("D12").Formula = "=SUM(D4:D11)"
$ws.Range("E12").Formula = "=SUM(E4:E11)"
$ws.Range("F12").Formula = "=SUM(F4:F11)"
$ws.Range("G12").Formula = "=SUM(G4:G11)"
$ws.Range("I12").Formula = "=SUM(I4:I11)"

# Row 13 ("Sınıf ORTALAMASI"): column averages
$ws.Range("D13").Formula = "=AVERAGE(D4:D12)"
$ws.Range("E13:G13").Formula = "=AVERAGE(E4:E12)"
$ws.Range("I13").Formula = "=SUM(I5:I12)"

# Student info block
$ws.Range("F17").Value = 20215070019
$ws.Range("F18").Value = "KÜBRA ÇABUK"
$ws.Range("F19").Value = "YBS"

# Update the selected range to match the saved view state
$ws.Range("F19:H19").Select()
